$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos crudos")

# Update existing rows 2-37 (timestamps in C, temperature in E) and add new rows 38-45
$timestamps = @(
    "2023-12-10 21:00:36",
    "2023-12-10 21:01:35",
    "2023-12-10 21:02:35",
    "2023-12-10 21:03:34",
    "2023-12-10 21:04:34",
    "2023-12-10 21:05:33",
    "2023-12-10 21:06:33",
    "2023-12-10 21:07:32",
    "2023-12-10 21:08:32",
    "2023-12-10 21:09:31",
    "2023-12-10 21:10:31",
    "2023-12-10 21:11:30",
    "2023-12-10 21:12:30",
    "2023-12-10 21:13:29",
    "2023-12-10 21:14:29",
    "2023-12-10 21:15:28",
    "2023-12-10 21:16:28",
    "2023-12-10 21:17:27",
    "2023-12-10 21:18:27",
    "2023-12-10 21:19:26",
    "2023-12-10 21:20:26",
    "2023-12-10 21:21:25",
    "2023-12-10 21:22:25",
    "2023-12-10 21:23:24",
    "2023-12-10 21:24:24",
    "2023-12-10 21:25:23",
    "2023-12-10 21:26:23",
    "2023-12-10 21:27:22",
    "2023-12-10 21:28:22",
    "2023-12-10 21:29:21",
    "2023-12-10 21:30:21",
    "2023-12-10 21:31:20",
    "2023-12-10 21:32:20",
    "2023-12-10 21:33:19",
    "2023-12-10 21:34:19",
    "2023-12-10 21:35:18",
    "2023-12-10 21:36:18",
    "2023-12-10 21:37:17",
    "2023-12-10 21:38:17",
    "2023-12-10 21:39:16",
    "2023-12-10 21:40:16",
    "2023-12-10 21:41:15",
    "2023-12-10 21:42:15",
    "2023-12-10 21:43:15"
)

$temps = @(
    25.017482517482499,
    22.0017482517482,
    20.034965034965001,
    18.986013986013901,
    17.412587412587399,
    16.363636363636299,
    15.576923076923,
    15.0524475524475,
    15.0524475524475,
    14.3968531468531,
    14.1346153846153,
    14.0034965034964,
    13.741258741258701,
    13.6101398601398,
    13.347902097902001,
    13.741258741258701,
    13.4790209790209,
    13.216783216783099,
    13.4790209790209,
    13.347902097902001,
    13.216783216783099,
    12.9545454545454,
    12.9545454545454,
    12.9545454545454,
    12.692307692307599,
    12.9545454545454,
    12.9545454545454,
    12.9545454545454,
    12.9545454545454,
    12.692307692307599,
    12.8234265734265,
    12.692307692307599,
    12.9545454545454,
    12.4300699300699,
    12.9545454545454,
    12.692307692307599,
    12.9545454545454,
    12.9545454545454,
    12.9545454545454,
    12.8234265734265,
    13.085664335664299,
    12.4300699300699,
    12.692307692307599,
    15.1835664335663
)


for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    if ($row -gt 37) {
        # New rows need the text number format applied (matches existing timestamp column style)
        $ws.Cells.Item($row, 3).NumberFormat = "@"
        $ws.Cells.Item($row, 1).Value = 6
        $ws.Cells.Item($row, 2).Value = 9
        $ws.Cells.Item($row, 4).Value = 0
    }
    $ws.Cells.Item($row, 3).Value = $timestamps[$i]
    $ws.Cells.Item($row, 5).Value = $temps[$i]
}

$ws.Range("H2").Value = 23
$ws.Range("H3").Formula = "=COUNT(E:E)-1"
